# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets to reflect refreshed numbers from the data source.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F5").Value = 882
    $ws.Range("F7").Value = 426
}
